$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell reference -> new text value.
# Values are written with a leading apostrophe so Excel always stores them
# as literal text (matching the original inlineStr cells), even for strings
# that look numeric (e.g. "1.00", "0.0000101"); ClearFormats() afterwards
# drops the quote-prefix cell style that the apostrophe would otherwise add,
# so number formatting stays identical to the untouched cells.
$updates = @(
    @{ Cell = 'D2'; Value = '46.431.45' }
    @{ Cell = 'E2'; Value = '  +1.04%  ' }
    @{ Cell = 'D3'; Value = '2.603.26' }
    @{ Cell = 'E3'; Value = '  +6.24%  ' }
    @{ Cell = 'E4'; Value = '  +0.19%  ' }
    @{ Cell = 'D5'; Value = '307.26' }
    @{ Cell = 'E5'; Value = '  +4.44%  ' }
    @{ Cell = 'D6'; Value = '100.51' }
    @{ Cell = 'E6'; Value = '  +5.43%  ' }
    @{ Cell = 'D7'; Value = '0.601' }
    @{ Cell = 'E7'; Value = '  +5.15%  ' }
    @{ Cell = 'E8'; Value = '  +0.09%  ' }
    @{ Cell = 'D9'; Value = '0.576' }
    @{ Cell = 'E9'; Value = '  +11.38%  ' }
    @{ Cell = 'D10'; Value = '39.49' }
    @{ Cell = 'E10'; Value = '  +12.99%  ' }
    @{ Cell = 'D11'; Value = '0.0846' }
    @{ Cell = 'E11'; Value = '  +8.03%  ' }
    @{ Cell = 'D12'; Value = '54.29' }
    @{ Cell = 'E12'; Value = '  +1.08%  ' }
    @{ Cell = 'D13'; Value = '8.16' }
    @{ Cell = 'E13'; Value = '  +12.36%  ' }
    @{ Cell = 'D14'; Value = '3.018.47' }
    @{ Cell = 'E14'; Value = '  +6.86%  ' }
    @{ Cell = 'E15'; Value = '  +1.97%  ' }
    @{ Cell = 'D16'; Value = '2.627.91' }
    @{ Cell = 'E16'; Value = '  +7.81%  ' }
    @{ Cell = 'D17'; Value = '0.925' }
    @{ Cell = 'E17'; Value = '  +9.00%  ' }
    @{ Cell = 'D18'; Value = '14.95' }
    @{ Cell = 'E18'; Value = '  +5.29%  ' }
    @{ Cell = 'D19'; Value = '46.599.54' }
    @{ Cell = 'E19'; Value = '  +1.48%  ' }
    @{ Cell = 'D20'; Value = '0.0000101' }
    @{ Cell = 'E20'; Value = '  +7.37%  ' }
    @{ Cell = 'D21'; Value = '12.93' }
    @{ Cell = 'E21'; Value = '  +2.10%  ' }
    @{ Cell = 'D22'; Value = '6.73' }
    @{ Cell = 'E22'; Value = '  +7.36%  ' }
    @{ Cell = 'B23'; Value = 'BitcoinCash' }
    @{ Cell = 'C23'; Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch' }
    @{ Cell = 'D23'; Value = '275.96' }
    @{ Cell = 'E23'; Value = '  +13.09%  ' }
    @{ Cell = 'B24'; Value = 'Litecoin' }
    @{ Cell = 'C24'; Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc' }
    @{ Cell = 'D24'; Value = '71.47' }
    @{ Cell = 'E24'; Value = '  +6.02%  ' }
    @{ Cell = 'D25'; Value = '3.03' }
    @{ Cell = 'E25'; Value = '  +9.12%  ' }
    @{ Cell = 'D26'; Value = '2.16' }
    @{ Cell = 'E26'; Value = '  +11.42%  ' }
    @{ Cell = 'D27'; Value = '28.88' }
    @{ Cell = 'E27'; Value = '  +33.93%  ' }
    @{ Cell = 'E28'; Value = '  -0.09%  ' }
    @{ Cell = 'D29'; Value = '4.03' }
    @{ Cell = 'E29'; Value = '  -0.44%  ' }
    @{ Cell = 'D30'; Value = '10.60' }
    @{ Cell = 'E30'; Value = '  +8.31%  ' }
    @{ Cell = 'B31'; Value = 'Toncoin' }
    @{ Cell = 'C31'; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton' }
    @{ Cell = 'D31'; Value = '2.31' }
    @{ Cell = 'E31'; Value = '  +4.22%  ' }
    @{ Cell = 'B32'; Value = 'InjectiveProtocol' }
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj' }
    @{ Cell = 'D32'; Value = '38.95' }
    @{ Cell = 'E32'; Value = '  +0.05%  ' }
    @{ Cell = 'D33'; Value = '6.34' }
    @{ Cell = 'E33'; Value = '  +14.53%  ' }
    @{ Cell = 'D34'; Value = '3.62' }
    @{ Cell = 'E34'; Value = '  -3.35%  ' }
    @{ Cell = 'E35'; Value = '  +4.65%  ' }
    @{ Cell = 'D36'; Value = '2.22' }
    @{ Cell = 'E36'; Value = '  +9.39%  ' }
    @{ Cell = 'D37'; Value = '0.0840' }
    @{ Cell = 'E37'; Value = '  +9.29%  ' }
    @{ Cell = 'D38'; Value = '151.37' }
    @{ Cell = 'E38'; Value = '  +2.92%  ' }
    @{ Cell = 'D39'; Value = '0.122' }
    @{ Cell = 'E39'; Value = '  +7.50%  ' }
    @{ Cell = 'E40'; Value = '  +6.50%  ' }
    @{ Cell = 'D41'; Value = '23.29' }
    @{ Cell = 'E41'; Value = '  +44.11%  ' }
    @{ Cell = 'D42'; Value = '15.92' }
    @{ Cell = 'E42'; Value = '  +5.97%  ' }
    @{ Cell = 'D43'; Value = '3.66' }
    @{ Cell = 'E43'; Value = '  +13.58%  ' }
    @{ Cell = 'D44'; Value = '0.0332' }
    @{ Cell = 'E44'; Value = '  +10.85%  ' }
    @{ Cell = 'D45'; Value = '4.07' }
    @{ Cell = 'E45'; Value = '  +3.67%  ' }
    @{ Cell = 'D46'; Value = '2.131.01' }
    @{ Cell = 'E46'; Value = '  +6.17%  ' }
    @{ Cell = 'D47'; Value = '1.00' }
    @{ Cell = 'E47'; Value = '  +0.17%  ' }
    @{ Cell = 'D48'; Value = '93.24' }
    @{ Cell = 'E48'; Value = '  +1.60%  ' }
    @{ Cell = 'D49'; Value = '9.51' }
    @{ Cell = 'E49'; Value = '  +11.29%  ' }
    @{ Cell = 'D50'; Value = '1.79' }
    @{ Cell = 'E50'; Value = '  +1.10%  ' }
    @{ Cell = 'D51'; Value = '108.94' }
    @{ Cell = 'E51'; Value = '  +6.45%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.Value = "'" + $u.Value
    $cell.ClearFormats()
}
